$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Hello"
$ws.Range("B7").Value = "World"
$ws.Range("C7").Value = "2025-10-01T18:34:25.432Z"
